$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows 4-7 correspond to the "Ready for handoff" entries.
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-21 04:39:20"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-21 04:39:24"

    $overview.Cells.Item($r, 7).Value = "2016-08-21 04:39:24"
}
